$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("A1").Value = "CS320 Section 101 (M-W-F 8:00 - 8:50)"
$ws.Range("A1:E1").Select()
